$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "GLAVNI TOK" (main flow) description text in B9: the old wording
# ("usposlenik ce iskoristiti sistem") is replaced with the new wording
# ("automatski ce se aktivirati sistem").
$ws.Range("B9").Value = "Nakon isteka roka prijava, automatski će se aktivirati sistem za zatvaranje dalje mogućnosti prijave za taj rok. Nakon toga će administrator sistema pokrenuti tombolu koja koristeći svoje sisteme odabira će odabrati primljene/odbijene vize."

# Adjust row heights to match the re-wrapped text (auto-fit results baked into the file).
$ws.Rows.Item(2).RowHeight = 63.85
$ws.Rows.Item(3).RowHeight = 60.45
$ws.Rows.Item(4).RowHeight = 63.2
$ws.Rows.Item(9).RowHeight = 76.1
$ws.Rows.Item(10).RowHeight = 97.15

# Update the active selection left behind after the edit.
$ws.Range("B9").Select()
